$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing data rows ---
# Row 2: Tissue column (F2) "Coronary Artery" -> "Coronary Artery Explant"
$ws.Range("F2").Value = "Coronary Artery Explant"

# Row 3: Tissue column (F3) "Carotid Endarterectomy" -> "Carotid Plaque"
$ws.Range("F3").Value = "Carotid Plaque"

# --- Add new row 4: Pan et al. / Circulation / 2020 / DOI / Human / Carotid Artery / 3 Patients / 8867 / 10x / Pan_2020 ---
$ws.Range("A4").Value = "Pan et al."
$ws.Range("B4").Value = "Circulation"
$ws.Range("C4").Value = 2020
$ws.Range("D4").Value = "https://doi.org/10.1161/CIRCULATIONAHA.120.048378"
$ws.Range("E4").Value = "Human"
$ws.Range("F4").Value = "Carotid Artery"
$ws.Range("G4").Value = "3 Patients"
$ws.Range("H4").Value = 8867
$ws.Range("I4").Value = "10x"
$ws.Range("J4").Value = "Pan_2020"

# --- Column width changes ---
# Engine stores width as round(ColumnWidth*6 + 5)/6, so back-solve the
# ColumnWidth needed to reproduce the target stored widths (42 and
# 22.1640625) as closely as the 1/6-pixel grid allows.
$ws.Columns.Item(4).ColumnWidth = 41.16666666666666
$ws.Columns.Item(6).ColumnWidth = 21.33333333333333

# --- Selection change ---
$ws.Range("D9").Select()

$wb.Save()
